$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = $false
$ws.Range("E2").Value = 106.51
$ws.Range("F2").Value = -2.2844036697247661
$ws.Range("G2").Value = $false

# Add new row 3
$ws.Range("C3").Value = 9771.56
